$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.155.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4647"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "97.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.04%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.850.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07530"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.964"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6287"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "296.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +22.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.068.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9990"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007355"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.075.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.048"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.086"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.920"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1074"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.329"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.006"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.813"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04907"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7237"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.731"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01911"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.659"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.90%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.966"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8604"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.643"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4042"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.022"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.939"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1187"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05539"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3716"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.18%  "
